$d = $word.ActiveDocument

# 1. Update the date string literal (search only the digits so the
#    surrounding straight quote characters are left untouched by
#    autocorrect/smart-quote substitution)
$d.Content.Find.Execute("2024-01-11", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-22", 2)

# 2. Remove the obsolete AIC.svrepglm code block (comments + function def + trailing blank lines)
$codeBlock = "# for our survey models, we want modelsummary() to retrieve the AIC value that is stored in the model object`v# instead of trying to calculate it as there's some sort of bug there (still to ask on StackOverflow about this)`vAIC.svrepglm <- function(model){`v  AIC <- model`$aic`v  `v  list(AIC = AIC)`v}`v`v`v"
$d.Content.Find.Execute($codeBlock, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3. Update the table numbers (unique values in the document)
$d.Content.Find.Execute("106.1", $true, $false, $false, $false, $false, $true, 1, $false, "-1925.7", 2)
$d.Content.Find.Execute("93.1", $true, $false, $false, $false, $false, $true, 1, $false, "-1934.0", 2)
$d.Content.Find.Execute("94.8", $true, $false, $false, $false, $false, $true, 1, $false, "-1930.3", 2)

$d.Content.Find.Execute("9226.2", $true, $false, $false, $false, $false, $true, 1, $false, "17496704.3", 2)
$d.Content.Find.Execute("9151.8", $true, $false, $false, $false, $false, $true, 1, $false, "17869013.4", 2)
$d.Content.Find.Execute("9105.0", $true, $false, $false, $false, $false, $true, 1, $false, "18515491.6", 2)

$d.Content.Find.Execute("-4602.064", $true, $false, $false, $false, $false, $true, 1, $false, "-8748341.109", 2)
$d.Content.Find.Execute("-4561.158", $true, $false, $false, $false, $false, $true, 1, $false, "-8934491.943", 2)
$d.Content.Find.Execute("-4534.057", $true, $false, $false, $false, $false, $true, 1, $false, "-9257727.361", 2)
